# PlayerPerformance_4654.xlsx edit:
#  - Add "Player Info" sheet (new first sheet)
#  - Keep "ODI Batting" and "ODI Bowling" sheets (now shifted to positions 2 and 3)
#    but change their MATCH_CARD_LINK column (a full URL) into a MATCH_CODE column
#    (just the numeric match code that used to be the MatchCode= query parameter)
#  - Add "ODI Batting Extra" sheet (new last sheet) with additional per-match
#    batting detail

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Add the new sheets in the right place so the final order is:
#    Player Info, ODI Batting, ODI Bowling, ODI Batting Extra
# ---------------------------------------------------------------------------
$firstSheet = $wb.Worksheets.Item(1)
$playerInfo = $wb.Worksheets.Add($firstSheet)
$playerInfo.Name = "Player Info"

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$battingExtra = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$battingExtra.Name = "ODI Batting Extra"

$battingSheet = $wb.Worksheets.Item("ODI Batting")
$bowlingSheet = $wb.Worksheets.Item("ODI Bowling")

# ---------------------------------------------------------------------------
# Helper used for header cells: bold + thin border + centered/top alignment
# (matches the workbook's existing header style)
# ---------------------------------------------------------------------------
function Format-HeaderCell($cell) {
    $cell.Font.Bold = $true
    $cell.Borders.LineStyle = 1
    $cell.HorizontalAlignment = -4108
    $cell.VerticalAlignment = -4160
}

# ---------------------------------------------------------------------------
# 2. Populate "Player Info"
# ---------------------------------------------------------------------------
$piHeaders = @("ID", "NAME", "BATTING_HAND", "BOWL_STYLE")
for ($c = 1; $c -le $piHeaders.Length; $c++) {
    $cell = $playerInfo.Cells.Item(1, $c)
    $cell.Value = $piHeaders[$c - 1]
    Format-HeaderCell $cell
}

# Column A holds a numeric-looking player id that must stay text, like the
# rest of this workbook's "numeric" columns.
$playerInfo.Range("A2").NumberFormat = "@"
$piRow = @("4654", "Fakhar Zaman", "Left Handed", "Left Arm Orthodox")
for ($c = 1; $c -le $piRow.Length; $c++) {
    $playerInfo.Cells.Item(2, $c).Value = $piRow[$c - 1]
}

# ---------------------------------------------------------------------------
# 3. "ODI Batting": rename MATCH_CARD_LINK -> MATCH_CODE (column D) and
#    replace each URL value with just the trailing MatchCode number (text)
# ---------------------------------------------------------------------------
$battingSheet.Cells.Item(1, 4).Value = "MATCH_CODE"
$battingLastRow = $battingSheet.UsedRange.Rows.Count
$battingSheet.Range($battingSheet.Cells.Item(2, 4), $battingSheet.Cells.Item($battingLastRow, 4)).NumberFormat = "@"
for ($r = 2; $r -le $battingLastRow; $r++) {
    $cell = $battingSheet.Cells.Item($r, 4)
    $v = $cell.Value()
    if ($v -ne $null) {
        $cell.Value = ($v -replace '.*MatchCode=', '')
    }
}

# ---------------------------------------------------------------------------
# 4. "ODI Bowling": rename MATCH_CARD_LINK -> MATCH_CODE (column B) and
#    replace each URL value with just the trailing MatchCode number (text)
# ---------------------------------------------------------------------------
$bowlingSheet.Cells.Item(1, 2).Value = "MATCH_CODE"
$bowlingLastRow = $bowlingSheet.UsedRange.Rows.Count
$bowlingSheet.Range($bowlingSheet.Cells.Item(2, 2), $bowlingSheet.Cells.Item($bowlingLastRow, 2)).NumberFormat = "@"
for ($r = 2; $r -le $bowlingLastRow; $r++) {
    $cell = $bowlingSheet.Cells.Item($r, 2)
    $v = $cell.Value()
    if ($v -ne $null) {
        $cell.Value = ($v -replace '.*MatchCode=', '')
    }
}

# ---------------------------------------------------------------------------
# 5. Populate "ODI Batting Extra"
# ---------------------------------------------------------------------------
$beHeaders = @("MATCH_CODE", "BATTING_POSITION", "NUM_4", "NUM_6", "PERCENT_RUNS_OF_TOTAL", "MAN_OF_MATCH")
for ($c = 1; $c -le $beHeaders.Length; $c++) {
    $cell = $battingExtra.Cells.Item(1, $c)
    $cell.Value = $beHeaders[$c - 1]
    Format-HeaderCell $cell
}

# MATCH_CODE (A), NUM_4 (C), NUM_6 (D) and PERCENT_RUNS_OF_TOTAL (E) are all
# numeric-looking but must stay text; BATTING_POSITION (B) is a real number.
$battingExtra.Range("A2:A21").NumberFormat = "@"

$beRows = @(
    @("4376", 1, "7", "1", "25.42%", "NO"),
    @("4434", 2, "0", "0", "0.72%", "NO"),
    @("4458", 2, "1", "0", "2.92%", "NO"),
    @("4459", 2, "18", "10", "59.57%", "YES"),
    @("4460", 2, "9", "3", "31.56%", "NO"),
    @("4472", 2, "6", "0", "33.33%", "NO"),
    @("4473", 2, "0", "0", "5.13%", "NO"),
    @("4476", 2, "1", "0", "1.81%", "NO"),
    @("4564", 1, "3", "0", "8.00%", "NO"),
    @("4565", $null, $null, $null, $null, "NO"),
    @("4567", 1, "3", "0", "7.94%", "NO"),
    @("4586", $null, $null, $null, $null, "NO"),
    @("4590", $null, $null, $null, $null, "NO"),
    @("4592", 1, "4", "0", "13.01%", "NO"),
    @("4634", $null, $null, $null, $null, "NO"),
    @("4638", 1, "0", "0", "1.57%", "NO"),
    @("4641", 2, "4", "0", "12.62%", "NO"),
    @("4686", $null, $null, $null, $null, "NO"),
    @("4688", 1, "0", "0", $null, "NO"),
    @("4690", $null, $null, $null, $null, "NO")
)

for ($i = 0; $i -lt $beRows.Length; $i++) {
    $row = $beRows[$i]
    $r = $i + 2
    for ($c = 1; $c -le $row.Length; $c++) {
        $val = $row[$c - 1]
        if ($val -ne $null) {
            $cell = $battingExtra.Cells.Item($r, $c)
            if ($c -eq 3 -or $c -eq 4 -or $c -eq 5) {
                $cell.NumberFormat = "@"
            }
            $cell.Value = $val
        }
    }
}
